{"js": "// Update the date line and every \"A\u00d7B=\" arithmetic cell per the commit's\n// regenerated worksheet content. Each original value is unique in the\n// document, so a plain exact-text search/replace is unambiguous for all of\n// them (including the header date line).\nconst replacements = [\n  [\"2024-09-19 Thursday\", \"2024-09-20 Friday\"],\n  [\"703\u00d72=\", \"989\u00d78=\"],\n  [\"908\u00d74=\", \"997\u00d77=\"],\n  [\"873\u00d77=\", \"528\u00d74=\"],\n  [\"888\u00d74=\", \"681\u00d75=\"],\n  [\"338\u00d73=\", \"240\u00d76=\"],\n  [\"430\u00d76=\", \"543\u00d75=\"],\n  [\"120\u00d72=\", \"947\u00d72=\"],\n  [\"795\u00d75=\", \"918\u00d75=\"],\n  [\"539\u00d73=\", \"408\u00d79=\"],\n  [\"101\u00d78=\", \"357\u00d78=\"],\n  [\"923\u00d73=\", \"929\u00d77=\"],\n  [\"609\u00d73=\", \"361\u00d72=\"],\n  [\"265\u00d76=\", \"356\u00d72=\"],\n  [\"938\u00d73=\", \"612\u00d79=\"],\n  [\"560\u00d73=\", \"319\u00d72=\"],\n  [\"129\u00d75=\", \"285\u00d75=\"],\n  [\"509\u00d74=\", \"635\u00d75=\"],\n  [\"549\u00d74=\", \"973\u00d77=\"],\n  [\"481\u00d78=\", \"413\u00d72=\"],\n  [\"356\u00d76=\", \"833\u00d77=\"],\n  [\"290\u00d73=\", \"583\u00d73=\"],\n  [\"456\u00d73=\", \"232\u00d76=\"],\n  [\"394\u00d77=\", \"121\u00d78=\"],\n  [\"196\u00d72=\", \"830\u00d74=\"],\n  [\"721\u00d73=\", \"861\u00d75=\"],\n];\n\nfor (const [oldText, newText] of replacements) {\n  const results = context.document.body.search(oldText, { matchCase: true, matchWholeWord: false });\n  results.load(\"items\");\n  await context.sync();\n\n  for (let i = 0; i < results.items.length; i++) {\n    results.items[i].insertText(newText, Word.InsertLocation.replace);\n  }\n  await context.sync();\n}\n", "ps1": "# Update the date line and every \"A\u00d7B=\" arithmetic cell per the commit's\n# regenerated worksheet content. Each original value is unique in the\n# document, so Find/Replace on the whole document Range is unambiguous for\n# every pair (including the header date line).\n$d = $word.ActiveDocument\n\n$replacements = @(\n    @(\"2024-09-19 Thursday\", \"2024-09-20 Friday\"),\n    @(\"703\u00d72=\", \"989\u00d78=\"),\n    @(\"908\u00d74=\", \"997\u00d77=\"),\n    @(\"873\u00d77=\", \"528\u00d74=\"),\n    @(\"888\u00d74=\", \"681\u00d75=\"),\n    @(\"338\u00d73=\", \"240\u00d76=\"),\n    @(\"430\u00d76=\", \"543\u00d75=\"),\n    @(\"120\u00d72=\", \"947\u00d72=\"),\n    @(\"795\u00d75=\", \"918\u00d75=\"),\n    @(\"539\u00d73=\", \"408\u00d79=\"),\n    @(\"101\u00d78=\", \"357\u00d78=\"),\n    @(\"923\u00d73=\", \"929\u00d77=\"),\n    @(\"609\u00d73=\", \"361\u00d72=\"),\n    @(\"265\u00d76=\", \"356\u00d72=\"),\n    @(\"938\u00d73=\", \"612\u00d79=\"),\n    @(\"560\u00d73=\", \"319\u00d72=\"),\n    @(\"129\u00d75=\", \"285\u00d75=\"),\n    @(\"509\u00d74=\", \"635\u00d75=\"),\n    @(\"549\u00d74=\", \"973\u00d77=\"),\n    @(\"481\u00d78=\", \"413\u00d72=\"),\n    @(\"356\u00d76=\", \"833\u00d77=\"),\n    @(\"290\u00d73=\", \"583\u00d73=\"),\n    @(\"456\u00d73=\", \"232\u00d76=\"),\n    @(\"394\u00d77=\", \"121\u00d78=\"),\n    @(\"196\u00d72=\", \"830\u00d74=\"),\n    @(\"721\u00d73=\", \"861\u00d75=\")\n)\n\nforeach ($pair in $replacements) {\n    $find = $pair[0]\n    $replace = $pair[1]\n    $rng = $d.Content\n    $rng.Find.Execute($find, $false, $false, $false, $false, $false, $true, 1, $false, $replace, 2)\n}\n"}
